$d = $word.ActiveDocument

# 1. "M.S. Computer Science – BS/MS Program " -> "M.S. Computer Science"
$d.Content.Find.Execute(
    "M.S. Computer Science " + [char]0x2013 + " BS/MS Program ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "M.S. Computer Science", 2) | Out-Null

# 2. "B.S. GPA – 3.65; M.S. GPA – 4.0" -> "B.S. GPA – 3.65, M.S. GPA – 4.0"
$d.Content.Find.Execute(
    "B.S. GPA " + [char]0x2013 + " 3.65; M.S. GPA " + [char]0x2013 + " 4.0",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "B.S. GPA " + [char]0x2013 + " 3.65, M.S. GPA " + [char]0x2013 + " 4.0", 2) | Out-Null

# 3. "Graduating with a B" -> "B" (remove "Graduating with a " prefix)
$d.Content.Find.Execute(
    "Graduating with a B",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "B", 2) | Out-Null

# 4. Remove the "Relevant Coursework: ..." bullet paragraph entirely, including
#    the following blank spacer paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("Relevant Coursework:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $rng.Paragraphs(1)
    $startRange = $para.Range
    # paragraph containing "Relevant Coursework" + the following blank paragraph
    $nextPara = $startRange.Next(4, 1).Paragraphs(1)
    $endOfBlank = $nextPara.Range.End
    $deleteRange = $d.Range($para.Range.Start, $endOfBlank)
    $deleteRange.Delete()
}

# 5. "NYC, New York (Remote) " -> "NYC, New York "
$d.Content.Find.Execute(
    "NYC, New York (Remote) ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "NYC, New York ", 2) | Out-Null

# 6. "Worked as a summer intern for the swaptions desk" ->
#    "Worked as a summer intern for the swap derivatives desk"
$d.Content.Find.Execute(
    "Worked as a summer intern for the swaptions desk",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Worked as a summer intern for the swap derivatives desk", 2) | Out-Null

# 7. "Computer Science 214: Data Structures and Data Management " ->
#    "Computer Science 214: Data Structures"
$d.Content.Find.Execute(
    "Computer Science 214: Data Structures and Data Management ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Computer Science 214: Data Structures", 2) | Out-Null
